$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1620.2433
$ws.Range("J17").Value = 1723.303
$ws.Range("L17").Value = 5169.909000000001
$ws.Range("N17").Value = -5505.909000000001
$ws.Range("H19").Value = 1466.8
$ws.Range("I19").Value = 401.25
$ws.Range("K19").Value = 401.25
$ws.Range("M19").Value = -226.25
$ws.Range("H131").Value = 3349.1
$ws.Range("I131").Value = 1912.7142
$ws.Range("J131").Value = 6700.6665
$ws.Range("K131").Value = 5738.142599999999
$ws.Range("L131").Value = 20101.9995
$ws.Range("M131").Value = -698.1425999999992
$ws.Range("N131").Value = -30181.9995
$ws.Range("H137").Value = 6325.3335
$ws.Range("I137").Value = 3954.2727
$ws.Range("J137").Value = 7955.4375
$ws.Range("K137").Value = 11862.8181
$ws.Range("L137").Value = 23866.3125
$ws.Range("M137").Value = -9312.8181
$ws.Range("N137").Value = -28966.3125
$ws.Range("H138").Value = 1453384.9
$ws.Range("J138").Value = 2637315.8
$ws.Range("L138").Value = 7911947.399999999
$ws.Range("N138").Value = -7922227.399999999
$ws.Range("H141").Value = 4030.8333
$ws.Range("I141").Value = 3284.75
$ws.Range("K141").Value = 9854.25
$ws.Range("M141").Value = -4674.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3339495
$ws.Range("I32").Value = 3395249
$ws.Range("K32").Value = 3395249
$ws.Range("M32").Value = -3394962
$ws.Range("H45").Value = 4407.316
$ws.Range("I45").Value = 1481.7
$ws.Range("J45").Value = 7658
$ws.Range("K45").Value = 1481.7
$ws.Range("L45").Value = 7658
$ws.Range("M45").Value = -1104.7
$ws.Range("N45").Value = -8412
$ws.Range("H97").Value = 6412584
$ws.Range("I97").Value = 2417.6667
$ws.Range("J97").Value = 20835458
$ws.Range("K97").Value = 2417.6667
$ws.Range("L97").Value = 20835458
$ws.Range("M97").Value = -1921.6667
$ws.Range("N97").Value = -20836450
$ws.Range("H122").Value = 3994.4
$ws.Range("I122").Value = 3250
$ws.Range("K122").Value = 9750
$ws.Range("M122").Value = -7300

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8340254
$ws.Range("I134").Value = 17859366
$ws.Range("K134").Value = 53578098
$ws.Range("M134").Value = -53575563

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 841.3333
$ws.Range("J22").Value = 350
$ws.Range("L22").Value = 350
$ws.Range("N22").Value = -1050
$ws.Range("H31").Value = 7156.022
$ws.Range("J31").Value = 10818.808
$ws.Range("L31").Value = 10818.808
$ws.Range("N31").Value = -11408.808
$ws.Range("H34").Value = 7156.022
$ws.Range("J34").Value = 10818.808
$ws.Range("L34").Value = 10818.808
$ws.Range("N34").Value = -11222.808
$ws.Range("H86").Value = 14209817
$ws.Range("I86").Value = 18387952
$ws.Range("J86").Value = 4157.8
$ws.Range("K86").Value = 18387952
$ws.Range("L86").Value = 4157.8
$ws.Range("M86").Value = -18386829
$ws.Range("N86").Value = -6403.8
$ws.Range("H89").Value = 14209817
$ws.Range("I89").Value = 18387952
$ws.Range("J89").Value = 4157.8
$ws.Range("K89").Value = 91939760
$ws.Range("L89").Value = 20789
$ws.Range("M89").Value = -91934144
$ws.Range("N89").Value = -32021
$ws.Range("H94").Value = 2646.5
$ws.Range("I94").Value = 2908
$ws.Range("K94").Value = 2908
$ws.Range("M94").Value = -2457
$ws.Range("H99").Value = 4889.706
$ws.Range("I99").Value = 2975.1428
$ws.Range("J99").Value = 6229.9
$ws.Range("K99").Value = 2975.1428
$ws.Range("L99").Value = 6229.9
$ws.Range("M99").Value = -1477.1428
$ws.Range("N99").Value = -9225.9
$ws.Range("H126").Value = 4889.706
$ws.Range("I126").Value = 2975.1428
$ws.Range("J126").Value = 6229.9
$ws.Range("K126").Value = 8925.428400000001
$ws.Range("L126").Value = 18689.7
$ws.Range("M126").Value = -6455.428400000001
$ws.Range("N126").Value = -23629.7
$ws.Range("H132").Value = 6612.6055
$ws.Range("I132").Value = 3989.85
$ws.Range("K132").Value = 11969.55
$ws.Range("M132").Value = -9439.549999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 48644476
$ws.Range("I4").Value = 54957384
$ws.Range("K4").Value = 164872152
$ws.Range("M4").Value = -164872040
$ws.Range("H113").Value = 5693.1665
$ws.Range("I113").Value = 1384.3334
$ws.Range("J113").Value = 10002
$ws.Range("K113").Value = 4153.0002
$ws.Range("L113").Value = 30006
$ws.Range("M113").Value = -1983.0002
$ws.Range("N113").Value = -34346

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1638.2778
$ws.Range("I97").Value = 820.7
$ws.Range("J97").Value = 2660.25
$ws.Range("K97").Value = 820.7
$ws.Range("L97").Value = 2660.25
$ws.Range("M97").Value = -324.7
$ws.Range("N97").Value = -3652.25
$ws.Range("H102").Value = 3052.6155
$ws.Range("I102").Value = 3279.3635
$ws.Range("K102").Value = 3279.3635
$ws.Range("M102").Value = -1657.3635
$ws.Range("H122").Value = 14329091
$ws.Range("I122").Value = 14329091
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 42987273
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -42984823
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4434.0586
$ws.Range("I132").Value = 1730.3478
$ws.Range("K132").Value = 5191.0434
$ws.Range("M132").Value = -2661.0434

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4032.0454
$ws.Range("I7").Value = 3100.2942
$ws.Range("K7").Value = 3100.2942
$ws.Range("M7").Value = -2988.2942
$ws.Range("H40").Value = 3101.838
$ws.Range("I40").Value = 2141.7856
$ws.Range("K40").Value = 2141.7856
$ws.Range("M40").Value = -2005.7856
$ws.Range("H82").Value = 521198.72
$ws.Range("I82").Value = 754786.75
$ws.Range("J82").Value = 2114.2222
$ws.Range("K82").Value = 754786.75
$ws.Range("L82").Value = 2114.2222
$ws.Range("M82").Value = -754425.75
$ws.Range("N82").Value = -2836.2222
$ws.Range("H85").Value = 521198.72
$ws.Range("I85").Value = 754786.75
$ws.Range("J85").Value = 2114.2222
$ws.Range("K85").Value = 754786.75
$ws.Range("L85").Value = 2114.2222
$ws.Range("M85").Value = -753538.75
$ws.Range("N85").Value = -4610.2222
$ws.Range("H93").Value = 1866.1666
$ws.Range("I93").Value = 1065.6666
$ws.Range("K93").Value = 1065.6666
$ws.Range("M93").Value = 182.3334
$ws.Range("H126").Value = 4032.0454
$ws.Range("I126").Value = 3100.2942
$ws.Range("K126").Value = 9300.882599999999
$ws.Range("M126").Value = -6830.882599999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10311.889
$ws.Range("I132").Value = 15401.5
$ws.Range("J132").Value = 6240.2
$ws.Range("K132").Value = 46204.5
$ws.Range("L132").Value = 18720.6
$ws.Range("M132").Value = -43674.5
$ws.Range("N132").Value = -23780.6
$ws.Range("H136").Value = 50506616
$ws.Range("I136").Value = 333334000
$ws.Range("J136").Value = 595903.5
$ws.Range("K136").Value = 1000002000
$ws.Range("L136").Value = 1787710.5
$ws.Range("M136").Value = -999999450
$ws.Range("N136").Value = -1792810.5
